$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; this shifts old rows 25,26,27 down to 26,27,28
$ws.Rows.Item(25).Insert()

# Row 24: label update only (data values unchanged)
$ws.Range("A24").Value = "Task 4.1"

# Row 25: "Task 4.2"
$ws.Range("A25").Value = "Task 4.2"
$ws.Range("B25").Value = 0.01057013667475462
$ws.Range("C25").Value = 0.0005005659407863555
$ws.Range("D25").Value = 0.01030731884798878
$ws.Range("E25").Value = 0.01046173435880525
$ws.Range("F25").Value = 0.01103629893751503
$ws.Range("G25").Value = 0.01008439654807944
$ws.Range("H25").Value = 0.01093962419329085
$ws.Range("I25").Value = 0.01015618973961673
$ws.Range("J25").Value = 0.0100724451293896
$ws.Range("K25").Value = 0.01061491848467104
$ws.Range("L25").Value = 0.01033790987375893
$ws.Range("M25").Value = 0.01048788691785228
$ws.Range("N25").Value = 0.01143696616118029
$ws.Range("O25").Value = 0.01092643645155179
$ws.Range("P25").Value = 0.01025622278850475
$ws.Range("Q25").Value = 0.01073533645900154
$ws.Range("R25").Value = 0.01007063401650483
$ws.Range("S25").Value = 0.01082403369716303
$ws.Range("T25").Value = 0.01050390026902077
$ws.Range("U25").Value = 0.01077713568897991
$ws.Range("V25").Value = 0.01156918574693464
$ws.Range("W25").Value = 0.009712423806234684
$ws.Range("X25").Value = 0.01169137996563375
$ws.Range("Y25").Value = 0.01008737208311145
$ws.Range("Z25").Value = 0.01038533426990044
$ws.Range("AA25").Value = 0.01002207990846811
$ws.Range("AB25").Value = 0.010892256726454
$ws.Range("AC25").Value = 0.01135329705292846
$ws.Range("AD25").Value = 0.01023900159137059
$ws.Range("AE25").Value = 0.01082215320353932
$ws.Range("AF25").Value = 0.01046445840502846
$ws.Range("AG25").Value = 0.009835768920159778

# Row 26: "Task 5.1"
$ws.Range("A26").Value = "Task 5.1"
$ws.Range("B26").Value = 0.01094645632520933
$ws.Range("C26").Value = 0.0009015060449485058
$ws.Range("D26").Value = 0.0105947442228127
$ws.Range("E26").Value = 0.01064116525558864
$ws.Range("F26").Value = 0.01290009441114279
$ws.Range("G26").Value = 0.01041402377965366
$ws.Range("H26").Value = 0.01052461469724162
$ws.Range("I26").Value = 0.01102616963801676
$ws.Range("J26").Value = 0.01043513092124959
$ws.Range("K26").Value = 0.01007429458308798
$ws.Range("L26").Value = 0.01052191798454783
$ws.Range("M26").Value = 0.01038607140137711
$ws.Range("N26").Value = 0.01169278225280685
$ws.Range("O26").Value = 0.01121431387448157
$ws.Range("P26").Value = 0.01013278970275092
$ws.Range("Q26").Value = 0.01054587546034712
$ws.Range("R26").Value = 0.009789165836331638
$ws.Range("S26").Value = 0.01057606753487713
$ws.Range("T26").Value = 0.01060314049240994
$ws.Range("U26").Value = 0.01088286083218916
$ws.Range("V26").Value = 0.01169587789401826
$ws.Range("W26").Value = 0.01036945468160288
$ws.Range("X26").Value = 0.01229077041653039
$ws.Range("Y26").Value = 0.01038786501902274
$ws.Range("Z26").Value = 0.01114223586665325
$ws.Range("AA26").Value = 0.01080742570073047
$ws.Range("AB26").Value = 0.0131134659350269
$ws.Range("AC26").Value = 0.01209163466560852
$ws.Range("AD26").Value = 0.009911310621027784
$ws.Range("AE26").Value = 0.01292687250175204
$ws.Range("AF26").Value = 0.01086647707732529
$ws.Range("AG26").Value = 0.009835076496068487

# Row 27: "Task 5.2"
$ws.Range("A27").Value = "Task 5.2"
$ws.Range("B27").Value = 0.01105812470007133
$ws.Range("C27").Value = 0.000801530499824401
$ws.Range("D27").Value = 0.01040134141374233
$ws.Range("E27").Value = 0.01054155830563586
$ws.Range("F27").Value = 0.01199250720794866
$ws.Range("G27").Value = 0.01017678899381777
$ws.Range("H27").Value = 0.01132816402016193
$ws.Range("I27").Value = 0.01194004847267183
$ws.Range("J27").Value = 0.01060357283730544
$ws.Range("K27").Value = 0.01131292985119929
$ws.Range("L27").Value = 0.01089105513727888
$ws.Range("M27").Value = 0.01026045188892836
$ws.Range("N27").Value = 0.01141995560590436
$ws.Range("O27").Value = 0.01168005161730688
$ws.Range("P27").Value = 0.010527887523696
$ws.Range("Q27").Value = 0.01036495060336992
$ws.Range("R27").Value = 0.009999550651592373
$ws.Range("S27").Value = 0.01054947398991609
$ws.Range("T27").Value = 0.01089902656314451
$ws.Range("U27").Value = 0.01152356853482261
$ws.Range("V27").Value = 0.01190753122639767
$ws.Range("W27").Value = 0.01074064948256935
$ws.Range("X27").Value = 0.012598873125015
$ws.Range("Y27").Value = 0.01076845084808825
$ws.Range("Z27").Value = 0.01174774694955943
$ws.Range("AA27").Value = 0.01032694570052849
$ws.Range("AB27").Value = 0.01168678052568842
$ws.Range("AC27").Value = 0.01152006667137058
$ws.Range("AD27").Value = 0.009958183749771426
$ws.Range("AE27").Value = 0.01325483482330768
$ws.Range("AF27").Value = 0.01089260112675032
$ws.Range("AG27").Value = 0.009928193554650223

# Row 28: "Task 5.3"
$ws.Range("A28").Value = "Task 5.3"
$ws.Range("B28").Value = 0.01092233396509257
$ws.Range("C28").Value = 0.0008218368162363563
$ws.Range("D28").Value = 0.01069095317834767
$ws.Range("E28").Value = 0.01046998944436819
$ws.Range("F28").Value = 0.01298152631821786
$ws.Range("G28").Value = 0.01037306942016312
$ws.Range("H28").Value = 0.0103084622219599
$ws.Range("I28").Value = 0.01100685503390813
$ws.Range("J28").Value = 0.01102118488917939
$ws.Range("K28").Value = 0.01112867270852016
$ws.Range("L28").Value = 0.01061221438384759
$ws.Range("M28").Value = 0.01046101830713468
$ws.Range("N28").Value = 0.01143602718992988
$ws.Range("O28").Value = 0.01110305415961764
$ws.Range("P28").Value = 0.01033951649990365
$ws.Range("Q28").Value = 0.01050756696355912
$ws.Range("R28").Value = 0.009873911427863166
$ws.Range("S28").Value = 0.01042816544272811
$ws.Range("T28").Value = 0.01096483406365935
$ws.Range("U28").Value = 0.01067821760227414
$ws.Range("V28").Value = 0.01258196408800248
$ws.Range("W28").Value = 0.01070875778470997
$ws.Range("X28").Value = 0.01247479091767482
$ws.Range("Y28").Value = 0.01055939637759182
$ws.Range("Z28").Value = 0.01082003997454337
$ws.Range("AA28").Value = 0.009938394700583848
$ws.Range("AB28").Value = 0.01265017724044342
$ws.Range("AC28").Value = 0.01082929094785161
$ws.Range("AD28").Value = 0.009938750174112372
$ws.Range("AE28").Value = 0.01205727783434322
$ws.Range("AF28").Value = 0.01075549715273997
$ws.Range("AG28").Value = 0.009970442504998421
